$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Row-level price / volume updates ---
Set-TextValue $ws.Range("D2") "51.884.37"
$ws.Range("E2").Value = "  -0.54%  "

Set-TextValue $ws.Range("D3") "2.783.66"
$ws.Range("E3").Value = "  -2.10%  "

Set-TextValue $ws.Range("D5") "356.91"
$ws.Range("E5").Value = "  -1.73%  "

Set-TextValue $ws.Range("D6") "109.30"
$ws.Range("E6").Value = "  -3.50%  "

$ws.Range("E7").Value = "  -2.58%  "

Set-TextValue $ws.Range("D8") "0.999"
$ws.Range("E8").Value = "  +0.00%  "

Set-TextValue $ws.Range("D9") "0.590"
$ws.Range("E9").Value = "  -2.27%  "

Set-TextValue $ws.Range("D10") "40.29"
$ws.Range("E10").Value = "  -3.48%  "

Set-TextValue $ws.Range("D11") "0.0848"
$ws.Range("E11").Value = "  -1.94%  "

$ws.Range("E12").Value = "  +1.36%  "

Set-TextValue $ws.Range("D13") "19.46"
$ws.Range("E13").Value = "  -3.33%  "

Set-TextValue $ws.Range("D14") "7.56"
$ws.Range("E14").Value = "  -3.25%  "

Set-TextValue $ws.Range("D15") "3.224.82"
$ws.Range("E15").Value = "  -1.92%  "

Set-TextValue $ws.Range("D16") "2.775.00"
$ws.Range("E16").Value = "  -2.16%  "

Set-TextValue $ws.Range("D17") "0.943"
$ws.Range("E17").Value = "  +2.75%  "

Set-TextValue $ws.Range("D18") "51.823.27"
$ws.Range("E18").Value = "  -0.42%  "

Set-TextValue $ws.Range("D19") "7.50"
$ws.Range("E19").Value = "  -1.01%  "

Set-TextValue $ws.Range("D20") "3.09"
$ws.Range("E20").Value = "  -2.49%  "

Set-TextValue $ws.Range("D21") "13.12"
$ws.Range("E21").Value = "  -3.68%  "

Set-TextValue $ws.Range("D22") "0.0₃0977"
$ws.Range("E22").Value = "  -2.39%  "

Set-TextValue $ws.Range("D23") "70.09"
$ws.Range("E23").Value = "  -0.45%  "

Set-TextValue $ws.Range("D24") "269.68"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("E25").Value = "  -3.83%  "

Set-TextValue $ws.Range("D26") "26.49"

$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("E28").Value = "  +15.81%  "

Set-TextValue $ws.Range("D29") "10.31"
$ws.Range("E29").Value = "  -0.94%  "

Set-TextValue $ws.Range("D30") "2.20"
$ws.Range("E30").Value = "  -2.24%  "

$ws.Range("E31").Value = "  -4.15%  "

Set-TextValue $ws.Range("D32") "52.11"
$ws.Range("E32").Value = "  -2.93%  "

Set-TextValue $ws.Range("D33") "34.40"
$ws.Range("E33").Value = "  -2.98%  "

Set-TextValue $ws.Range("D34") "5.73"
$ws.Range("E34").Value = "  -2.50%  "

Set-TextValue $ws.Range("D35") "0.0845"
$ws.Range("E35").Value = "  -0.18%  "

Set-TextValue $ws.Range("D36") "5.20"
$ws.Range("E36").Value = "  -4.69%  "

$ws.Range("E37").Value = "  +0.02%  "

Set-TextValue $ws.Range("D38") "18.91"
$ws.Range("E38").Value = "  +2.23%  "

Set-TextValue $ws.Range("D39") "3.20"
$ws.Range("E39").Value = "  -2.86%  "

$ws.Range("E40").Value = "  -4.44%  "

$ws.Range("E41").Value = "  +2.97%  "

Set-TextValue $ws.Range("D42") "0.115"
$ws.Range("E42").Value = "  -2.39%  "

$ws.Range("E43").Value = "  -1.07%  "

Set-TextValue $ws.Range("D44") "119.86"
$ws.Range("E44").Value = "  -5.45%  "

Set-TextValue $ws.Range("D45") "21.83"
$ws.Range("E45").Value = "  -8.07%  "

Set-TextValue $ws.Range("D46") "2.087.53"
$ws.Range("E46").Value = "  -1.39%  "

Set-TextValue $ws.Range("D47") "3.26"
$ws.Range("E47").Value = "  -5.00%  "

$ws.Range("E48").Value = "  -1.98%  "

Set-TextValue $ws.Range("D51") "1.14"
$ws.Range("E51").Value = "  +31.39%  "

# --- Rows 49 and 50 swapped (THORChain <-> SEI) plus updated price/volume ---
$ws.Range("B49").Value = "SEI"
$ws.Range("C49").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
Set-TextValue $ws.Range("D49") "0.960"
$ws.Range("E49").Value = "  -2.85%  "

$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D50") "5.74"
$ws.Range("E50").Value = "  -2.61%  "

